$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Cells.Item(16, 2).Value = 6528871
$ws.Cells.Item(16, 5).Value = 'VSK Arhus'
$ws.Cells.Item(16, 6).Value = 'Ishoj'
$ws.Cells.Item(16, 7).Value = 1
$ws.Cells.Item(16, 8).Value = 3
$ws.Cells.Item(16, 9).Value = 1
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 'A'
$ws.Cells.Item(16, 12).Value = 1.6
$ws.Cells.Item(16, 13).Value = 4
$ws.Cells.Item(16, 14).Value = 4.5
$ws.Cells.Item(16, 15).Value = 1.3
$ws.Cells.Item(16, 16).Value = 5
$ws.Cells.Item(16, 17).Value = 7
$ws.Cells.Item(16, 18).Value = -1.75
$ws.Cells.Item(16, 19).Value = 1.95
$ws.Cells.Item(16, 20).Value = 1.85
$ws.Cells.Item(16, 21).Value = 3.5
$ws.Cells.Item(16, 22).Value = 1.95
$ws.Cells.Item(16, 23).Value = 1.85
$ws.Cells.Item(16, 24).Value = -1
$ws.Cells.Item(16, 25).Value = -1
$ws.Cells.Item(16, 26).Value = 6
$ws.Cells.Item(16, 27).Value = -1
$ws.Cells.Item(16, 28).Value = 0.8500000000000001
$ws.Cells.Item(16, 29).Value = 0.95
$ws.Cells.Item(16, 30).Value = -1

# Row 17
$ws.Cells.Item(17, 2).Value = 6526660
$ws.Cells.Item(17, 5).Value = 'IF Lyseng'
$ws.Cells.Item(17, 6).Value = 'Young Boys FD'
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 1
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 'A'
$ws.Cells.Item(17, 12).Value = 3
$ws.Cells.Item(17, 13).Value = 3.6
$ws.Cells.Item(17, 14).Value = 2
$ws.Cells.Item(17, 15).Value = 3.1
$ws.Cells.Item(17, 16).Value = 4
$ws.Cells.Item(17, 17).Value = 1.85
$ws.Cells.Item(17, 18).Value = 0.25
$ws.Cells.Item(17, 19).Value = 2.1
$ws.Cells.Item(17, 20).Value = 1.7
$ws.Cells.Item(17, 21).Value = 3.25
$ws.Cells.Item(17, 22).Value = 1.7
$ws.Cells.Item(17, 23).Value = 2.1
$ws.Cells.Item(17, 24).Value = -1
$ws.Cells.Item(17, 25).Value = -1
$ws.Cells.Item(17, 26).Value = 0.8500000000000001
$ws.Cells.Item(17, 27).Value = -1
$ws.Cells.Item(17, 28).Value = 0.7
$ws.Cells.Item(17, 29).Value = -1
$ws.Cells.Item(17, 30).Value = 1.1

# Row 18
$ws.Cells.Item(18, 2).Value = 6526519
$ws.Cells.Item(18, 5).Value = 'Holstebro'
$ws.Cells.Item(18, 6).Value = 'AB Trnby'
$ws.Cells.Item(18, 7).Value = 4
$ws.Cells.Item(18, 8).Value = 2
$ws.Cells.Item(18, 9).Value = 2
$ws.Cells.Item(18, 10).Value = 2
$ws.Cells.Item(18, 11).Value = 'H'
$ws.Cells.Item(18, 12).Value = 1.571
$ws.Cells.Item(18, 13).Value = 4.2
$ws.Cells.Item(18, 14).Value = 4.333
$ws.Cells.Item(18, 15).Value = 1.5
$ws.Cells.Item(18, 16).Value = 4.5
$ws.Cells.Item(18, 17).Value = 5
$ws.Cells.Item(18, 18).Value = -1
$ws.Cells.Item(18, 19).Value = 1.8
$ws.Cells.Item(18, 20).Value = 2
$ws.Cells.Item(18, 21).Value = 3.25
$ws.Cells.Item(18, 22).Value = 1.8
$ws.Cells.Item(18, 23).Value = 2
$ws.Cells.Item(18, 24).Value = 0.5
$ws.Cells.Item(18, 25).Value = -1
$ws.Cells.Item(18, 26).Value = -1
$ws.Cells.Item(18, 27).Value = 0.8
$ws.Cells.Item(18, 28).Value = -1
$ws.Cells.Item(18, 29).Value = 0.8
$ws.Cells.Item(18, 30).Value = -1

# Row 19
$ws.Cells.Item(19, 2).Value = 6528884
$ws.Cells.Item(19, 5).Value = 'Nsby'
$ws.Cells.Item(19, 6).Value = 'FA 2000'
$ws.Cells.Item(19, 7).Value = 4
$ws.Cells.Item(19, 8).Value = 3
$ws.Cells.Item(19, 9).Value = 1
$ws.Cells.Item(19, 10).Value = 2
$ws.Cells.Item(19, 11).Value = 'H'
$ws.Cells.Item(19, 12).Value = 2.875
$ws.Cells.Item(19, 13).Value = 3.8
$ws.Cells.Item(19, 14).Value = 2.05
$ws.Cells.Item(19, 15).Value = 2.8
$ws.Cells.Item(19, 16).Value = 3.75
$ws.Cells.Item(19, 17).Value = 2.1
$ws.Cells.Item(19, 18).Value = 0.25
$ws.Cells.Item(19, 19).Value = 1.875
$ws.Cells.Item(19, 20).Value = 1.925
$ws.Cells.Item(19, 21).Value = 3
$ws.Cells.Item(19, 22).Value = 2
$ws.Cells.Item(19, 23).Value = 1.8
$ws.Cells.Item(19, 24).Value = 1.8
$ws.Cells.Item(19, 25).Value = -1
$ws.Cells.Item(19, 26).Value = -1
$ws.Cells.Item(19, 27).Value = 0.875
$ws.Cells.Item(19, 28).Value = -1
$ws.Cells.Item(19, 29).Value = 1
$ws.Cells.Item(19, 30).Value = -1

# Row 22
$ws.Cells.Item(22, 2).Value = 6858502
$ws.Cells.Item(22, 5).Value = 'IF Lyseng'
$ws.Cells.Item(22, 6).Value = 'Ishoj'
$ws.Cells.Item(22, 7).Value = 2
$ws.Cells.Item(22, 8).Value = 2
$ws.Cells.Item(22, 9).Value = 1
$ws.Cells.Item(22, 10).Value = 1
$ws.Cells.Item(22, 11).Value = 'D'
$ws.Cells.Item(22, 12).Value = 2.875
$ws.Cells.Item(22, 13).Value = 3.75
$ws.Cells.Item(22, 14).Value = 2
$ws.Cells.Item(22, 15).Value = 3.4
$ws.Cells.Item(22, 16).Value = 4
$ws.Cells.Item(22, 17).Value = 1.75
$ws.Cells.Item(22, 18).Value = 0.5
$ws.Cells.Item(22, 19).Value = 1.975
$ws.Cells.Item(22, 20).Value = 1.825
$ws.Cells.Item(22, 21).Value = 3.25
$ws.Cells.Item(22, 22).Value = 1.85
$ws.Cells.Item(22, 23).Value = 1.95
$ws.Cells.Item(22, 24).Value = -1
$ws.Cells.Item(22, 25).Value = 3
$ws.Cells.Item(22, 26).Value = -1
$ws.Cells.Item(22, 27).Value = 0.9750000000000001
$ws.Cells.Item(22, 28).Value = -1
$ws.Cells.Item(22, 29).Value = 0.8500000000000001
$ws.Cells.Item(22, 30).Value = -1

# Row 23
$ws.Cells.Item(23, 2).Value = 6859123
$ws.Cells.Item(23, 5).Value = 'Vanlse'
$ws.Cells.Item(23, 6).Value = 'Avarta'
$ws.Cells.Item(23, 7).Value = 3
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 2
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 'H'
$ws.Cells.Item(23, 12).Value = 1.666
$ws.Cells.Item(23, 13).Value = 3.75
$ws.Cells.Item(23, 14).Value = 4
$ws.Cells.Item(23, 15).Value = 1.85
$ws.Cells.Item(23, 16).Value = 3.8
$ws.Cells.Item(23, 17).Value = 3.25
$ws.Cells.Item(23, 18).Value = -0.5
$ws.Cells.Item(23, 19).Value = 1.9
$ws.Cells.Item(23, 20).Value = 1.9
$ws.Cells.Item(23, 21).Value = 3
$ws.Cells.Item(23, 22).Value = 1.8
$ws.Cells.Item(23, 23).Value = 2
$ws.Cells.Item(23, 24).Value = 0.8500000000000001
$ws.Cells.Item(23, 25).Value = -1
$ws.Cells.Item(23, 26).Value = -1
$ws.Cells.Item(23, 27).Value = 0.8999999999999999
$ws.Cells.Item(23, 28).Value = -1
$ws.Cells.Item(23, 29).Value = 0
$ws.Cells.Item(23, 30).Value = 0

# Row 52
$ws.Cells.Item(52, 2).Value = 6858784
$ws.Cells.Item(52, 5).Value = 'Vejgaard B'
$ws.Cells.Item(52, 6).Value = 'Oure FA'
$ws.Cells.Item(52, 7).Value = 4
$ws.Cells.Item(52, 8).Value = 1
$ws.Cells.Item(52, 9).Value = 1
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 11).Value = 'H'
$ws.Cells.Item(52, 12).Value = 1.65
$ws.Cells.Item(52, 13).Value = 3.75
$ws.Cells.Item(52, 14).Value = 4.333
$ws.Cells.Item(52, 15).Value = 1.571
$ws.Cells.Item(52, 16).Value = 3.8
$ws.Cells.Item(52, 17).Value = 4.75
$ws.Cells.Item(52, 18).Value = -1
$ws.Cells.Item(52, 19).Value = 1.925
$ws.Cells.Item(52, 20).Value = 1.875
$ws.Cells.Item(52, 21).Value = 3.25
$ws.Cells.Item(52, 22).Value = 1.9
$ws.Cells.Item(52, 23).Value = 1.9
$ws.Cells.Item(52, 24).Value = 0.571
$ws.Cells.Item(52, 25).Value = -1
$ws.Cells.Item(52, 26).Value = -1
$ws.Cells.Item(52, 27).Value = 0.925
$ws.Cells.Item(52, 28).Value = -1
$ws.Cells.Item(52, 29).Value = 0.8999999999999999
$ws.Cells.Item(52, 30).Value = -1

# Row 53
$ws.Cells.Item(53, 2).Value = 6858783
$ws.Cells.Item(53, 5).Value = 'IF Lyseng'
$ws.Cells.Item(53, 6).Value = 'Vanlse'
$ws.Cells.Item(53, 7).Value = 4
$ws.Cells.Item(53, 8).Value = 2
$ws.Cells.Item(53, 9).Value = 2
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 'H'
$ws.Cells.Item(53, 12).Value = 3.25
$ws.Cells.Item(53, 13).Value = 3.75
$ws.Cells.Item(53, 14).Value = 1.85
$ws.Cells.Item(53, 15).Value = 3.2
$ws.Cells.Item(53, 16).Value = 3.5
$ws.Cells.Item(53, 17).Value = 2
$ws.Cells.Item(53, 18).Value = 0.25
$ws.Cells.Item(53, 19).Value = 1.975
$ws.Cells.Item(53, 20).Value = 1.825
$ws.Cells.Item(53, 21).Value = 2.75
$ws.Cells.Item(53, 22).Value = 1.925
$ws.Cells.Item(53, 23).Value = 1.875
$ws.Cells.Item(53, 24).Value = 2.2
$ws.Cells.Item(53, 25).Value = -1
$ws.Cells.Item(53, 26).Value = -1
$ws.Cells.Item(53, 27).Value = 0.9750000000000001
$ws.Cells.Item(53, 28).Value = -1
$ws.Cells.Item(53, 29).Value = 0.925
$ws.Cells.Item(53, 30).Value = -1

# Row 60
$ws.Cells.Item(60, 2).Value = 6858790
$ws.Cells.Item(60, 5).Value = 'Nsby'
$ws.Cells.Item(60, 6).Value = 'IF Lyseng'
$ws.Cells.Item(60, 7).Value = 2
$ws.Cells.Item(60, 8).Value = 2
$ws.Cells.Item(60, 9).Value = 2
$ws.Cells.Item(60, 10).Value = 1
$ws.Cells.Item(60, 11).Value = 'D'
$ws.Cells.Item(60, 12).Value = 2.35
$ws.Cells.Item(60, 13).Value = 3.6
$ws.Cells.Item(60, 14).Value = 2.5
$ws.Cells.Item(60, 15).Value = 2.2
$ws.Cells.Item(60, 16).Value = 3.75
$ws.Cells.Item(60, 17).Value = 2.625
$ws.Cells.Item(60, 18).Value = -0.25
$ws.Cells.Item(60, 19).Value = 1.975
$ws.Cells.Item(60, 20).Value = 1.825
$ws.Cells.Item(60, 21).Value = 3.25
$ws.Cells.Item(60, 22).Value = 1.9
$ws.Cells.Item(60, 23).Value = 1.9
$ws.Cells.Item(60, 24).Value = -1
$ws.Cells.Item(60, 25).Value = 2.75
$ws.Cells.Item(60, 26).Value = -1
$ws.Cells.Item(60, 27).Value = -0.5
$ws.Cells.Item(60, 28).Value = 0.4125
$ws.Cells.Item(60, 29).Value = 0.8999999999999999
$ws.Cells.Item(60, 30).Value = -1

# Row 61
$ws.Cells.Item(61, 2).Value = 6859110
$ws.Cells.Item(61, 5).Value = 'Young Boys FD'
$ws.Cells.Item(61, 6).Value = 'Avarta'
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 2
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 1
$ws.Cells.Item(61, 11).Value = 'A'
$ws.Cells.Item(61, 12).Value = 1.5
$ws.Cells.Item(61, 13).Value = 4.2
$ws.Cells.Item(61, 14).Value = 5.25
$ws.Cells.Item(61, 15).Value = 1.5
$ws.Cells.Item(61, 16).Value = 4.2
$ws.Cells.Item(61, 17).Value = 5.5
$ws.Cells.Item(61, 18).Value = -1
$ws.Cells.Item(61, 19).Value = 1.8
$ws.Cells.Item(61, 20).Value = 2
$ws.Cells.Item(61, 21).Value = 3
$ws.Cells.Item(61, 22).Value = 1.95
$ws.Cells.Item(61, 23).Value = 1.85
$ws.Cells.Item(61, 24).Value = -1
$ws.Cells.Item(61, 25).Value = -1
$ws.Cells.Item(61, 26).Value = 4.5
$ws.Cells.Item(61, 27).Value = -1
$ws.Cells.Item(61, 28).Value = 1
$ws.Cells.Item(61, 29).Value = -1
$ws.Cells.Item(61, 30).Value = 0.8500000000000001

# Row 62
$ws.Cells.Item(62, 2).Value = 6858791
$ws.Cells.Item(62, 5).Value = 'IF Lyseng'
$ws.Cells.Item(62, 6).Value = 'Young Boys FD'
$ws.Cells.Item(62, 7).Value = 3
$ws.Cells.Item(62, 8).Value = 3
$ws.Cells.Item(62, 9).Value = 1
$ws.Cells.Item(62, 10).Value = 3
$ws.Cells.Item(62, 11).Value = 'D'
$ws.Cells.Item(62, 12).Value = 2.8
$ws.Cells.Item(62, 13).Value = 3.6
$ws.Cells.Item(62, 14).Value = 2.15
$ws.Cells.Item(62, 15).Value = 2.9
$ws.Cells.Item(62, 16).Value = 3.4
$ws.Cells.Item(62, 17).Value = 2.15
$ws.Cells.Item(62, 18).Value = 0.25
$ws.Cells.Item(62, 19).Value = 1.875
$ws.Cells.Item(62, 20).Value = 1.925
$ws.Cells.Item(62, 21).Value = 3
$ws.Cells.Item(62, 22).Value = 1.875
$ws.Cells.Item(62, 23).Value = 1.925
$ws.Cells.Item(62, 24).Value = -1
$ws.Cells.Item(62, 25).Value = 2.4
$ws.Cells.Item(62, 26).Value = -1
$ws.Cells.Item(62, 27).Value = 0.4375
$ws.Cells.Item(62, 28).Value = -0.5
$ws.Cells.Item(62, 29).Value = 0.875
$ws.Cells.Item(62, 30).Value = -1

# Row 63
$ws.Cells.Item(63, 2).Value = 6858792
$ws.Cells.Item(63, 5).Value = 'Ishoj'
$ws.Cells.Item(63, 6).Value = 'Oure FA'
$ws.Cells.Item(63, 7).Value = 2
$ws.Cells.Item(63, 8).Value = 2
$ws.Cells.Item(63, 9).Value = 2
$ws.Cells.Item(63, 10).Value = 2
$ws.Cells.Item(63, 11).Value = 'D'
$ws.Cells.Item(63, 12).Value = 1.45
$ws.Cells.Item(63, 13).Value = 4.5
$ws.Cells.Item(63, 14).Value = 5.5
$ws.Cells.Item(63, 15).Value = 1.45
$ws.Cells.Item(63, 16).Value = 4.5
$ws.Cells.Item(63, 17).Value = 5.25
$ws.Cells.Item(63, 18).Value = -1.25
$ws.Cells.Item(63, 19).Value = 1.925
$ws.Cells.Item(63, 20).Value = 1.875
$ws.Cells.Item(63, 21).Value = 3.25
$ws.Cells.Item(63, 22).Value = 1.875
$ws.Cells.Item(63, 23).Value = 1.925
$ws.Cells.Item(63, 24).Value = -1
$ws.Cells.Item(63, 25).Value = 3.5
$ws.Cells.Item(63, 26).Value = -1
$ws.Cells.Item(63, 27).Value = -1
$ws.Cells.Item(63, 28).Value = 0.875
$ws.Cells.Item(63, 29).Value = 0.875
$ws.Cells.Item(63, 30).Value = -1

# Row 76
$ws.Cells.Item(76, 2).Value = 6858800
$ws.Cells.Item(76, 5).Value = 'Frem'
$ws.Cells.Item(76, 6).Value = 'Young Boys FD'
$ws.Cells.Item(76, 7).Value = 1
$ws.Cells.Item(76, 8).Value = 2
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 1
$ws.Cells.Item(76, 11).Value = 'A'
$ws.Cells.Item(76, 12).Value = 2.3
$ws.Cells.Item(76, 13).Value = 3.5
$ws.Cells.Item(76, 14).Value = 2.625
$ws.Cells.Item(76, 15).Value = 2.3
$ws.Cells.Item(76, 16).Value = 3.5
$ws.Cells.Item(76, 17).Value = 2.625
$ws.Cells.Item(76, 18).Value = -0.25
$ws.Cells.Item(76, 19).Value = 2.05
$ws.Cells.Item(76, 20).Value = 1.75
$ws.Cells.Item(76, 21).Value = 3
$ws.Cells.Item(76, 22).Value = 1.9
$ws.Cells.Item(76, 23).Value = 1.9
$ws.Cells.Item(76, 24).Value = -1
$ws.Cells.Item(76, 25).Value = -1
$ws.Cells.Item(76, 26).Value = 1.625
$ws.Cells.Item(76, 27).Value = -1
$ws.Cells.Item(76, 28).Value = 0.75
$ws.Cells.Item(76, 29).Value = 0
$ws.Cells.Item(76, 30).Value = 0

# Row 77
$ws.Cells.Item(77, 2).Value = 6858801
$ws.Cells.Item(77, 5).Value = 'Nsby'
$ws.Cells.Item(77, 6).Value = 'Vanlse'
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(77, 8).Value = 1
$ws.Cells.Item(77, 9).Value = 1
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 'D'
$ws.Cells.Item(77, 12).Value = 2.45
$ws.Cells.Item(77, 13).Value = 3.6
$ws.Cells.Item(77, 14).Value = 2.5
$ws.Cells.Item(77, 15).Value = 2.45
$ws.Cells.Item(77, 16).Value = 3.6
$ws.Cells.Item(77, 17).Value = 2.5
$ws.Cells.Item(77, 18).Value = 0
$ws.Cells.Item(77, 19).Value = 1.875
$ws.Cells.Item(77, 20).Value = 1.925
$ws.Cells.Item(77, 21).Value = 3
$ws.Cells.Item(77, 22).Value = 1.975
$ws.Cells.Item(77, 23).Value = 1.825
$ws.Cells.Item(77, 24).Value = -1
$ws.Cells.Item(77, 25).Value = 2.6
$ws.Cells.Item(77, 26).Value = -1
$ws.Cells.Item(77, 27).Value = 0
$ws.Cells.Item(77, 28).Value = 0
$ws.Cells.Item(77, 29).Value = -1
$ws.Cells.Item(77, 30).Value = 0.825

# Row 99
$ws.Cells.Item(99, 2).Value = 6858815
$ws.Cells.Item(99, 5).Value = 'IF Lyseng'
$ws.Cells.Item(99, 6).Value = 'Nsby'
$ws.Cells.Item(99, 7).Value = 1
$ws.Cells.Item(99, 8).Value = 1
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 1
$ws.Cells.Item(99, 11).Value = 'D'
$ws.Cells.Item(99, 12).Value = 2.1
$ws.Cells.Item(99, 13).Value = 3.5
$ws.Cells.Item(99, 14).Value = 3
$ws.Cells.Item(99, 15).Value = 2.4
$ws.Cells.Item(99, 16).Value = 3.2
$ws.Cells.Item(99, 17).Value = 2.7
$ws.Cells.Item(99, 18).Value = 0
$ws.Cells.Item(99, 19).Value = 1.775
$ws.Cells.Item(99, 20).Value = 2.025
$ws.Cells.Item(99, 21).Value = 2.75
$ws.Cells.Item(99, 22).Value = 1.975
$ws.Cells.Item(99, 23).Value = 1.825
$ws.Cells.Item(99, 24).Value = -1
$ws.Cells.Item(99, 25).Value = 2.2
$ws.Cells.Item(99, 26).Value = -1
$ws.Cells.Item(99, 27).Value = 0
$ws.Cells.Item(99, 28).Value = 0
$ws.Cells.Item(99, 29).Value = -1
$ws.Cells.Item(99, 30).Value = 0.825

# Row 100
$ws.Cells.Item(100, 2).Value = 6858818
$ws.Cells.Item(100, 5).Value = 'Vanlse'
$ws.Cells.Item(100, 6).Value = 'Oure FA'
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 3
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = 1
$ws.Cells.Item(100, 11).Value = 'A'
$ws.Cells.Item(100, 12).Value = 1.65
$ws.Cells.Item(100, 13).Value = 3.8
$ws.Cells.Item(100, 14).Value = 4.5
$ws.Cells.Item(100, 15).Value = 1.533
$ws.Cells.Item(100, 16).Value = 4.333
$ws.Cells.Item(100, 17).Value = 4.75
$ws.Cells.Item(100, 18).Value = -1
$ws.Cells.Item(100, 19).Value = 1.825
$ws.Cells.Item(100, 20).Value = 1.975
$ws.Cells.Item(100, 21).Value = 3.25
$ws.Cells.Item(100, 22).Value = 1.9
$ws.Cells.Item(100, 23).Value = 1.9
$ws.Cells.Item(100, 24).Value = -1
$ws.Cells.Item(100, 25).Value = -1
$ws.Cells.Item(100, 26).Value = 3.75
$ws.Cells.Item(100, 27).Value = -1
$ws.Cells.Item(100, 28).Value = 0.9750000000000001
$ws.Cells.Item(100, 29).Value = -0.5
$ws.Cells.Item(100, 30).Value = 0.45

# Row 101
$ws.Cells.Item(101, 2).Value = 6858817
$ws.Cells.Item(101, 5).Value = 'Vejgaard B'
$ws.Cells.Item(101, 6).Value = 'Young Boys FD'
$ws.Cells.Item(101, 7).Value = 1
$ws.Cells.Item(101, 8).Value = 2
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 1
$ws.Cells.Item(101, 11).Value = 'A'
$ws.Cells.Item(101, 12).Value = 2.45
$ws.Cells.Item(101, 13).Value = 3.75
$ws.Cells.Item(101, 14).Value = 2.4
$ws.Cells.Item(101, 15).Value = 2.15
$ws.Cells.Item(101, 16).Value = 3.8
$ws.Cells.Item(101, 17).Value = 2.75
$ws.Cells.Item(101, 18).Value = -0.25
$ws.Cells.Item(101, 19).Value = 1.95
$ws.Cells.Item(101, 20).Value = 1.85
$ws.Cells.Item(101, 21).Value = 3
$ws.Cells.Item(101, 22).Value = 1.75
$ws.Cells.Item(101, 23).Value = 1.95
$ws.Cells.Item(101, 24).Value = -1
$ws.Cells.Item(101, 25).Value = -1
$ws.Cells.Item(101, 26).Value = 1.75
$ws.Cells.Item(101, 27).Value = -1
$ws.Cells.Item(101, 28).Value = 0.8500000000000001
$ws.Cells.Item(101, 29).Value = 0
$ws.Cells.Item(101, 30).Value = 0

# Row 146
$ws.Cells.Item(146, 2).Value = 6858850
$ws.Cells.Item(146, 5).Value = 'Vejgaard B'
$ws.Cells.Item(146, 6).Value = 'Holstebro'
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 5
$ws.Cells.Item(146, 9).Value = 0
$ws.Cells.Item(146, 10).Value = 1
$ws.Cells.Item(146, 11).Value = 'A'
$ws.Cells.Item(146, 12).Value = 1.85
$ws.Cells.Item(146, 13).Value = 3.5
$ws.Cells.Item(146, 14).Value = 3.6
$ws.Cells.Item(146, 15).Value = 1.85
$ws.Cells.Item(146, 16).Value = 3.5
$ws.Cells.Item(146, 17).Value = 3.6
$ws.Cells.Item(146, 18).Value = -0.5
$ws.Cells.Item(146, 19).Value = 1.9
$ws.Cells.Item(146, 20).Value = 1.9
$ws.Cells.Item(146, 21).Value = 3
$ws.Cells.Item(146, 22).Value = 2
$ws.Cells.Item(146, 23).Value = 1.8
$ws.Cells.Item(146, 24).Value = -1
$ws.Cells.Item(146, 25).Value = -1
$ws.Cells.Item(146, 26).Value = 2.6
$ws.Cells.Item(146, 27).Value = -1
$ws.Cells.Item(146, 28).Value = 0.8999999999999999
$ws.Cells.Item(146, 29).Value = 1
$ws.Cells.Item(146, 30).Value = -1

# Row 147
$ws.Cells.Item(147, 2).Value = 6858849
$ws.Cells.Item(147, 5).Value = 'Frem'
$ws.Cells.Item(147, 6).Value = 'VSK Arhus'
$ws.Cells.Item(147, 7).Value = 2
$ws.Cells.Item(147, 8).Value = 0
$ws.Cells.Item(147, 9).Value = 0
$ws.Cells.Item(147, 10).Value = 0
$ws.Cells.Item(147, 11).Value = 'H'
$ws.Cells.Item(147, 12).Value = 2.05
$ws.Cells.Item(147, 13).Value = 3.4
$ws.Cells.Item(147, 14).Value = 3.1
$ws.Cells.Item(147, 15).Value = 1.833
$ws.Cells.Item(147, 16).Value = 3.5
$ws.Cells.Item(147, 17).Value = 3.6
$ws.Cells.Item(147, 18).Value = -0.5
$ws.Cells.Item(147, 19).Value = 1.875
$ws.Cells.Item(147, 20).Value = 1.925
$ws.Cells.Item(147, 21).Value = 2.5
$ws.Cells.Item(147, 22).Value = 1.85
$ws.Cells.Item(147, 23).Value = 1.95
$ws.Cells.Item(147, 24).Value = 0.833
$ws.Cells.Item(147, 25).Value = -1
$ws.Cells.Item(147, 26).Value = -1
$ws.Cells.Item(147, 27).Value = 0.875
$ws.Cells.Item(147, 28).Value = -1
$ws.Cells.Item(147, 29).Value = -1
$ws.Cells.Item(147, 30).Value = 0.95

# Row 148
$ws.Cells.Item(148, 2).Value = 6859081
$ws.Cells.Item(148, 5).Value = 'Holbk'
$ws.Cells.Item(148, 6).Value = 'Vanlse'
$ws.Cells.Item(148, 7).Value = 1
$ws.Cells.Item(148, 8).Value = 2
$ws.Cells.Item(148, 9).Value = 0
$ws.Cells.Item(148, 10).Value = 0
$ws.Cells.Item(148, 11).Value = 'A'
$ws.Cells.Item(148, 12).Value = 1.8
$ws.Cells.Item(148, 13).Value = 3.5
$ws.Cells.Item(148, 14).Value = 3.8
$ws.Cells.Item(148, 15).Value = 1.8
$ws.Cells.Item(148, 16).Value = 3.5
$ws.Cells.Item(148, 17).Value = 3.8
$ws.Cells.Item(148, 18).Value = -0.5
$ws.Cells.Item(148, 19).Value = 1.825
$ws.Cells.Item(148, 20).Value = 1.975
$ws.Cells.Item(148, 21).Value = 2.5
$ws.Cells.Item(148, 22).Value = 2
$ws.Cells.Item(148, 23).Value = 1.8
$ws.Cells.Item(148, 24).Value = -1
$ws.Cells.Item(148, 25).Value = -1
$ws.Cells.Item(148, 26).Value = 2.8
$ws.Cells.Item(148, 27).Value = -1
$ws.Cells.Item(148, 28).Value = 0.9750000000000001
$ws.Cells.Item(148, 29).Value = 1
$ws.Cells.Item(148, 30).Value = -1

# Row 164
$ws.Cells.Item(164, 2).Value = 8137133
$ws.Cells.Item(164, 5).Value = 'Ishoj'
$ws.Cells.Item(164, 6).Value = 'Frem'
$ws.Cells.Item(164, 7).Value = 1
$ws.Cells.Item(164, 8).Value = 1
$ws.Cells.Item(164, 9).Value = 1
$ws.Cells.Item(164, 10).Value = 0
$ws.Cells.Item(164, 11).Value = 'D'
$ws.Cells.Item(164, 12).Value = 3.1
$ws.Cells.Item(164, 13).Value = 3.5
$ws.Cells.Item(164, 14).Value = 2.05
$ws.Cells.Item(164, 15).Value = 3.1
$ws.Cells.Item(164, 16).Value = 3.75
$ws.Cells.Item(164, 17).Value = 2
$ws.Cells.Item(164, 18).Value = 0.5
$ws.Cells.Item(164, 19).Value = 1.75
$ws.Cells.Item(164, 20).Value = 1.95
$ws.Cells.Item(164, 21).Value = 3
$ws.Cells.Item(164, 22).Value = 1.9
$ws.Cells.Item(164, 23).Value = 1.9
$ws.Cells.Item(164, 24).Value = -1
$ws.Cells.Item(164, 25).Value = 2.75
$ws.Cells.Item(164, 26).Value = -1
$ws.Cells.Item(164, 27).Value = 0.75
$ws.Cells.Item(164, 28).Value = -1
$ws.Cells.Item(164, 29).Value = -1
$ws.Cells.Item(164, 30).Value = 0.8999999999999999

# Row 165
$ws.Cells.Item(165, 2).Value = 8089036
$ws.Cells.Item(165, 5).Value = 'Vejgaard B'
$ws.Cells.Item(165, 6).Value = 'IF Lyseng'
$ws.Cells.Item(165, 7).Value = 5
$ws.Cells.Item(165, 8).Value = 5
$ws.Cells.Item(165, 9).Value = 1
$ws.Cells.Item(165, 10).Value = 3
$ws.Cells.Item(165, 11).Value = 'D'
$ws.Cells.Item(165, 12).Value = 2.375
$ws.Cells.Item(165, 13).Value = 3.4
$ws.Cells.Item(165, 14).Value = 2.6
$ws.Cells.Item(165, 15).Value = 2.3
$ws.Cells.Item(165, 16).Value = 3.75
$ws.Cells.Item(165, 17).Value = 2.6
$ws.Cells.Item(165, 18).Value = 0
$ws.Cells.Item(165, 19).Value = 1.8
$ws.Cells.Item(165, 20).Value = 2
$ws.Cells.Item(165, 21).Value = 2.75
$ws.Cells.Item(165, 22).Value = 1.875
$ws.Cells.Item(165, 23).Value = 1.925
$ws.Cells.Item(165, 24).Value = -1
$ws.Cells.Item(165, 25).Value = 2.75
$ws.Cells.Item(165, 26).Value = -1
$ws.Cells.Item(165, 27).Value = 0
$ws.Cells.Item(165, 28).Value = 0
$ws.Cells.Item(165, 29).Value = 0.875
$ws.Cells.Item(165, 30).Value = -1
